$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tastes & aromas of varieties")
$ws.Activate()

$varieties = @(
    "Cabernet-Sauvignon",
    "Merlot",
    "Nebbiolo",
    "Pinot Noir",
    "Sangiovese",
    "Syrah/Shiraz",
    "Zinfandel",
    "Aglianico",
    "Gamay",
    "Barbera",
    "Cabernet franc",
    "Corvina veronese",
    "Grenache",
    "Malbec",
    "Mourvèdre",
    "Nerello Mascalese",
    "Nero d’Avola",
    "Tempranillo"
)

$startRow = 23
for ($i = 0; $i -lt $varieties.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $varieties[$i]
}

$ws.Range("A23").WrapText = $true

$ws.Range("A40").Select()
